$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H32").Value = 45460720
$ws.Range("I32").Value = 100007700
$ws.Range("K32").Value = 100007700
$ws.Range("M32").Value = -100007374
$ws.Range("H93").Value = 97499
$ws.Range("J93").Value = 97499
$ws.Range("L93").Value = 97499
$ws.Range("N93").Value = -102491
$ws.Range("H103").Value = 1315.3125
$ws.Range("J103").Value = 2377
$ws.Range("L103").Value = 7131
$ws.Range("N103").Value = -8303
$ws.Range("H132").Value = 32264066
$ws.Range("I132").Value = 37042236
$ws.Range("K132").Value = 111126708
$ws.Range("M132").Value = -111124178
$ws.Range("H135").Value = 11158.444
$ws.Range("I135").Value = 8575.143
$ws.Range("K135").Value = 77176.287
$ws.Range("M135").Value = -74641.287
$ws.Range("H137").Value = 3835.2766
$ws.Range("I137").Value = 3794.7307
$ws.Range("J137").Value = 3885.476
$ws.Range("K137").Value = 11384.1921
$ws.Range("L137").Value = 11656.428
$ws.Range("M137").Value = -8834.1921
$ws.Range("N137").Value = -16756.428
$ws.Range("H138").Value = 530399.5
$ws.Range("J138").Value = 718667.4
$ws.Range("L138").Value = 2156002.2
$ws.Range("N138").Value = -2166282.2
$ws.Range("H141").Value = 2976.8462
$ws.Range("I141").Value = 2976.8462
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 8930.5386
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -3750.5386
$ws.Range("N141").ClearContents()

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 1001.0941
$ws.Range("I32").Value = 946.22974
$ws.Range("K32").Value = 946.22974
$ws.Range("M32").Value = -659.22974
$ws.Range("H45").Value = 4047.353
$ws.Range("I45").Value = 3700.1304
$ws.Range("K45").Value = 3700.1304
$ws.Range("M45").Value = -3323.1304
$ws.Range("H61").Value = 6358.7915
$ws.Range("I61").Value = 7427.8335
$ws.Range("K61").Value = 7427.8335
$ws.Range("M61").Value = -7215.8335
$ws.Range("H74").Value = 4674.9443
$ws.Range("I74").Value = 5589.4
$ws.Range("K74").Value = 5589.4
$ws.Range("M74").Value = -4715.4
$ws.Range("H77").Value = 4674.9443
$ws.Range("I77").Value = 5589.4
$ws.Range("K77").Value = 27947
$ws.Range("M77").Value = -23579
$ws.Range("H122").Value = 1258.1666
$ws.Range("I122").Value = 1264.2307
$ws.Range("J122").Value = 1242.4
$ws.Range("K122").Value = 3792.6921
$ws.Range("L122").Value = 3727.2
$ws.Range("M122").Value = -1342.6921
$ws.Range("N122").Value = -8627.200000000001
$ws.Range("H132").Value = 4254.6787
$ws.Range("I132").Value = 2725.4614
$ws.Range("J132").Value = 4717
$ws.Range("K132").Value = 8176.3842
$ws.Range("L132").Value = 14151
$ws.Range("M132").Value = -5646.3842
$ws.Range("N132").Value = -19211
$ws.Range("H136").Value = 6358.7915
$ws.Range("I136").Value = 7427.8335
$ws.Range("K136").Value = 22283.5005
$ws.Range("M136").Value = -19733.5005

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H36").Value = 1362.6666
$ws.Range("I36").Value = 635.2
$ws.Range("K36").Value = 635.2
$ws.Range("M36").Value = -101.2
$ws.Range("H105").Value = 4159.407
$ws.Range("I105").Value = 2957.75
$ws.Range("K105").Value = 2957.75
$ws.Range("M105").Value = -1210.75

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 2655.4
$ws.Range("I31").Value = 1197.7941
$ws.Range("J31").Value = 3406.2878
$ws.Range("K31").Value = 1197.7941
$ws.Range("L31").Value = 3406.2878
$ws.Range("M31").Value = -902.7941000000001
$ws.Range("N31").Value = -3996.2878
$ws.Range("H34").Value = 2655.4
$ws.Range("I34").Value = 1197.7941
$ws.Range("J34").Value = 3406.2878
$ws.Range("K34").Value = 1197.7941
$ws.Range("L34").Value = 3406.2878
$ws.Range("M34").Value = -995.7941000000001
$ws.Range("N34").Value = -3810.2878
$ws.Range("H99").Value = 3969.261
$ws.Range("I99").Value = 4107.0527
$ws.Range("J99").Value = 3314.75
$ws.Range("K99").Value = 4107.0527
$ws.Range("L99").Value = 3314.75
$ws.Range("M99").Value = -2609.0527
$ws.Range("N99").Value = -6310.75
$ws.Range("H126").Value = 3969.261
$ws.Range("I126").Value = 4107.0527
$ws.Range("J126").Value = 3314.75
$ws.Range("K126").Value = 12321.1581
$ws.Range("L126").Value = 9944.25
$ws.Range("M126").Value = -9851.158100000001
$ws.Range("N126").Value = -14884.25
$ws.Range("H132").Value = 3481
$ws.Range("I132").Value = 2940.375
$ws.Range("K132").Value = 8821.125
$ws.Range("M132").Value = -6291.125
$ws.Range("H134").Value = 3698.2444
$ws.Range("I134").Value = 3555.9768
$ws.Range("K134").Value = 10667.9304
$ws.Range("M134").Value = -8132.930399999999

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H11").Value = 40
$ws.Range("I11").Value = 40
$ws.Range("K11").Value = 120
$ws.Range("M11").Value = 20
$ws.Range("H46").Value = 313049.9
$ws.Range("I46").Value = 999666.3
$ws.Range("J46").Value = 18785.715
$ws.Range("K46").Value = 2998998.9
$ws.Range("L46").Value = 56357.145
$ws.Range("M46").Value = -2998907.9
$ws.Range("N46").Value = -56539.145
$ws.Range("H107").Value = 838.7778
$ws.Range("J107").Value = 1149.8182
$ws.Range("L107").Value = 3449.4546
$ws.Range("N107").Value = -7289.4546
$ws.Range("H132").Value = 1602.7142
$ws.Range("I132").Value = 500.5
$ws.Range("J132").Value = 2043.6
$ws.Range("K132").Value = 4504.5
$ws.Range("L132").Value = 18392.4
$ws.Range("M132").Value = -1974.5
$ws.Range("N132").Value = -23452.4
$ws.Range("H137").Value = 25367580
$ws.Range("I137").Value = 25003206
$ws.Range("K137").Value = 75009618
$ws.Range("M137").Value = -75004518

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H96").Value = 57557.75
$ws.Range("J96").Value = 69998
$ws.Range("L96").Value = 69998
$ws.Range("N96").Value = -75490
$ws.Range("H102").Value = 15855.815
$ws.Range("I102").Value = 1832.56
$ws.Range("K102").Value = 1832.56
$ws.Range("M102").Value = -210.5599999999999
$ws.Range("H122").Value = 30366216
$ws.Range("I122").Value = 38541010
$ws.Range("J122").Value = 2698.4285
$ws.Range("K122").Value = 115623030
$ws.Range("L122").Value = 8095.2855
$ws.Range("M122").Value = -115620580
$ws.Range("N122").Value = -12995.2855
$ws.Range("H126").Value = 10794
$ws.Range("I126").Value = 6632.6665
$ws.Range("K126").Value = 19897.9995
$ws.Range("M126").Value = -17427.9995
$ws.Range("H132").Value = 5262.9033
$ws.Range("I132").Value = 8840.846
$ws.Range("K132").Value = 26522.538
$ws.Range("M132").Value = -23992.538

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H16").Value = 1942.92
$ws.Range("I16").Value = 817.2
$ws.Range("K16").Value = 817.2
$ws.Range("M16").Value = -647.2
$ws.Range("H61").Value = 24194.467
$ws.Range("I61").Value = 26522.12
$ws.Range("K61").Value = 26522.12
$ws.Range("M61").Value = -26320.12
$ws.Range("H113").Value = 24194.467
$ws.Range("I113").Value = 26522.12
$ws.Range("K113").Value = 26522.12
$ws.Range("M113").Value = -24352.12
$ws.Range("H132").Value = 9115.936
$ws.Range("I132").Value = 9087.559999999999
$ws.Range("J132").Value = 9234.166999999999
$ws.Range("K132").Value = 27262.68
$ws.Range("L132").Value = 27702.501
$ws.Range("M132").Value = -24732.68
$ws.Range("N132").Value = -32762.501
$ws.Range("H136").Value = 4074.7026
$ws.Range("I136").Value = 3893.5862
$ws.Range("K136").Value = 11680.7586
$ws.Range("M136").Value = -9130.758600000001

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H2").Value = 251005000
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H3").Value = 3348334
$ws.Range("I3").Value = 3348334
$ws.Range("K3").Value = 3348334
$ws.Range("M3").Value = -3348220
$ws.Range("H45").Value = 26371.5
$ws.Range("I45").Value = 22895
$ws.Range("K45").Value = 22895
$ws.Range("M45").Value = -22404
$ws.Range("H46").Value = 193333
$ws.Range("J46").Value = 193333
$ws.Range("L46").Value = 193333
$ws.Range("N46").Value = -193795
$ws.Range("H88").Value = 55355
$ws.Range("J88").Value = 55355
$ws.Range("L88").Value = 55355
$ws.Range("N88").Value = -56167
$ws.Range("H91").Value = 55355
$ws.Range("J91").Value = 55355
$ws.Range("L91").Value = 55355
$ws.Range("N91").Value = -58163
$ws.Range("H132").Value = 3008.6428
$ws.Range("I132").Value = 2312
$ws.Range("K132").Value = 6936
$ws.Range("M132").Value = -4406
$ws.Range("H134").Value = 193333
$ws.Range("J134").Value = 193333
$ws.Range("L134").Value = 579999
$ws.Range("N134").Value = -585069
$ws.Range("H136").Value = 4090
$ws.Range("I136").Value = 2188.682
$ws.Range("K136").Value = 6566.045999999999
$ws.Range("M136").Value = -4016.045999999999
